$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.75
$ws.Range("I3").Value = 4
$ws.Range("K3").Value = 21
$ws.Range("P3").Value = 1.25
$ws.Range("Q3").Value = 3.75
$ws.Range("R3").Value = 1.5
$ws.Range("S3").Value = 2.5
$ws.Range("T3").Value = 11
$ws.Range("U3").Value = 11
$ws.Range("X3").Value = 12
$ws.Range("AF3").Value = 26

# Row 4
$ws.Range("J4").Value = 1.05
$ws.Range("K4").Value = 11

# Row 10
$ws.Range("G10").Value = 2.05
$ws.Range("I10").Value = 3.8
$ws.Range("J10").Value = 1.11
$ws.Range("K10").Value = 6.5
$ws.Range("R10").Value = 2.5
$ws.Range("S10").Value = 1.5
$ws.Range("X10").Value = 23

# Row 12
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 3
$ws.Range("I12").Value = 2.37
$ws.Range("L12").Value = 1.47
$ws.Range("M12").Value = 2.32
$ws.Range("N12").Value = 2.35
$ws.Range("O12").Value = 1.47
$ws.Range("P12").Value = 1.52
$ws.Range("Q12").Value = 2.22
$ws.Range("R12").Value = 2
$ws.Range("T12").Value = 7.2
$ws.Range("U12").Value = 14
$ws.Range("V12").Value = 11.5
$ws.Range("W12").Value = 40
$ws.Range("X12").Value = 32
$ws.Range("Y12").Value = 50
$ws.Range("Z12").Value = 6.7
$ws.Range("AA12").Value = 5.9
$ws.Range("AB12").Value = 18
$ws.Range("AC12").Value = 110
$ws.Range("AE12").Value = 6.1
$ws.Range("AF12").Value = 10.25
$ws.Range("AG12").Value = 10
$ws.Range("AH12").Value = 24
$ws.Range("AI12").Value = 24

# Row 13
$ws.Range("P13").Value = 1.47
$ws.Range("Q13").Value = 2.32
$ws.Range("AJ13").Value = 45

# Row 15
$ws.Range("G15").Value = 3.95
$ws.Range("H15").Value = 2.75
$ws.Range("I15").Value = 2.18
$ws.Range("K15").Value = 4.9
$ws.Range("Q15").Value = 2.2
$ws.Range("R15").Value = 2.15
$ws.Range("S15").Value = 1.62
$ws.Range("T15").Value = 8
$ws.Range("U15").Value = 19.5
$ws.Range("X15").Value = 45
$ws.Range("Y15").Value = 65
$ws.Range("Z15").Value = 4.9
$ws.Range("AA15").Value = 5.5
$ws.Range("AB15").Value = 18.5
$ws.Range("AC15").Value = 120
$ws.Range("AE15").Value = 5.2
$ws.Range("AI15").Value = 23

# Row 16
$ws.Range("H16").Value = 3.15
$ws.Range("I16").Value = 1.78
$ws.Range("V16").Value = 16
$ws.Range("AB16").Value = 18
$ws.Range("AF16").Value = 7.1
$ws.Range("AH16").Value = 14
$ws.Range("AI16").Value = 17

# Row 17
$ws.Range("H17").Value = 2.95
$ws.Range("I17").Value = 2.5
$ws.Range("J17").Value = 1.08
$ws.Range("K17").Value = 6.4
$ws.Range("L17").Value = 1.37
$ws.Range("M17").Value = 2.87
$ws.Range("P17").Value = 1.4
$ws.Range("Q17").Value = 2.7
$ws.Range("T17").Value = 8.25
$ws.Range("V17").Value = 10.25
$ws.Range("X17").Value = 26
$ws.Range("Y17").Value = 35
$ws.Range("Z17").Value = 6.4
$ws.Range("AD17").Value = 500
$ws.Range("AE17").Value = 7.6
$ws.Range("AF17").Value = 12
$ws.Range("AH17").Value = 28
$ws.Range("AI17").Value = 22
$ws.Range("AJ17").Value = 32

# Row 18
$ws.Range("R18").Value = 1.73
$ws.Range("S18").Value = 2
$ws.Range("U18").Value = 11
$ws.Range("V18").Value = 9
$ws.Range("X18").Value = 17
$ws.Range("Y18").Value = 26
$ws.Range("AB18").Value = 13
$ws.Range("AC18").Value = 41
$ws.Range("AD18").Value = 201
$ws.Range("AE18").Value = 9.5
$ws.Range("AH18").Value = 34
$ws.Range("AI18").Value = 26
$ws.Range("AJ18").Value = 29

# Row 23
$ws.Range("G23").Value = 3.4
$ws.Range("H23").Value = 2.88
$ws.Range("I23").Value = 2.3
$ws.Range("W23").Value = 34
$ws.Range("Z23").Value = 7.5
$ws.Range("AF23").Value = 11
$ws.Range("AI23").Value = 21

# Row 26
$ws.Range("K26").Value = 13
$ws.Range("N26").Value = 1.75
$ws.Range("O26").Value = 2.05

# Row 27
$ws.Range("K27").Value = 13

# Row 28
$ws.Range("N28").Value = 1.83
$ws.Range("O28").Value = 1.98

# Row 33
$ws.Range("G33").Value = 1.83
$ws.Range("H33").Value = 3.4
$ws.Range("I33").Value = 4.35
$ws.Range("J33").Value = 1.09
$ws.Range("K33").Value = 6.6
$ws.Range("L33").Value = 1.42
$ws.Range("M33").Value = 2.72
$ws.Range("N33").Value = 2.25
$ws.Range("O33").Value = 1.6
$ws.Range("P33").Value = 1.52
$ws.Range("Q33").Value = 2.45
$ws.Range("R33").Value = 2.05
$ws.Range("U33").Value = 8
$ws.Range("V33").Value = 9.25
$ws.Range("W33").Value = 15.5
$ws.Range("X33").Value = 18
$ws.Range("Y33").Value = 40
$ws.Range("Z33").Value = 6.6
$ws.Range("AA33").Value = 7
$ws.Range("AB33").Value = 21
$ws.Range("AE33").Value = 9.5
$ws.Range("AF33").Value = 24
$ws.Range("AG33").Value = 16
$ws.Range("AH33").Value = 80
$ws.Range("AI33").Value = 55
$ws.Range("AJ33").Value = 70

# Row 34
$ws.Range("G34").Value = 4
$ws.Range("H34").Value = 3.4
$ws.Range("I34").Value = 1.95
$ws.Range("K34").Value = 8.5
$ws.Range("AB34").Value = 17
$ws.Range("AD34").Value = 351
$ws.Range("AE34").Value = 6.5
$ws.Range("AF34").Value = 8.5

# Row 41
$ws.Range("G41").Value = 3
$ws.Range("H41").Value = 3.55
$ws.Range("I41").Value = 2.15
$ws.Range("J41").Value = 1.03
$ws.Range("K41").Value = 9
$ws.Range("L41").Value = 1.18
$ws.Range("M41").Value = 4.25
$ws.Range("N41").Value = 1.57
$ws.Range("O41").Value = 2.27
$ws.Range("P41").Value = 1.3
$ws.Range("Q41").Value = 3.2
$ws.Range("T41").Value = 13
$ws.Range("U41").Value = 18.5
$ws.Range("V41").Value = 10.5
$ws.Range("W41").Value = 37
$ws.Range("X41").Value = 22
$ws.Range("Y41").Value = 24
$ws.Range("Z41").Value = 9
$ws.Range("AA41").Value = 7.3
$ws.Range("AB41").Value = 11.25
$ws.Range("AE41").Value = 11
$ws.Range("AF41").Value = 13
$ws.Range("AG41").Value = 8.75
$ws.Range("AH41").Value = 23
$ws.Range("AI41").Value = 15
$ws.Range("AJ41").Value = 19.5

# Row 42
$ws.Range("G42").Value = 3.2
$ws.Range("H42").Value = 3.6
$ws.Range("I42").Value = 2.05
$ws.Range("K42").Value = 8.5
$ws.Range("O42").Value = 2.18
$ws.Range("P42").Value = 1.32
$ws.Range("Q42").Value = 3.1
$ws.Range("R42").Value = 1.53
$ws.Range("S42").Value = 2.35
$ws.Range("U42").Value = 19.5
$ws.Range("V42").Value = 11
$ws.Range("X42").Value = 24
$ws.Range("Y42").Value = 26
$ws.Range("Z42").Value = 8.5
$ws.Range("AB42").Value = 11.75
$ws.Range("AC42").Value = 40
$ws.Range("AF42").Value = 11.75
$ws.Range("AH42").Value = 20
$ws.Range("AI42").Value = 15
